# Atualiza pontuacoes e resultados das competicoes
# Updates column V (Pontuacao) values on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 68.06
    3  = 53.66
    5  = 70
    6  = 58.96
    7  = 57.45
    8  = 68.06
    9  = 47.86
    10 = 64.56
    11 = 43.46
    12 = 61.56
    13 = 84.86
    14 = 42.96
    15 = 74.06
    16 = 54.6
    17 = 84.26000000000001
    18 = 59.25
    19 = 54.16
    20 = 63.9
    21 = 50.26
    22 = 73.95999999999999
    23 = 72.7
    24 = 57.6
    25 = 64.7
    27 = 49
    28 = 61.96
    29 = 60.16
    30 = 62.56
    31 = 73.76000000000001
    32 = 64.56
    33 = 47.86
}

foreach ($row in $updates.Keys) {
    $ws.Range("V$row").Value = $updates[$row]
}
